$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.896.43'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '2.569.83'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''302.46'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').Value = '''96.17'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.548'
$ws.Range('E9').Value = '  -1.85%  '
$ws.Range('D10').Value = '''36.32'
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').Value = '''0.0809'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '''7.62'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('E13').Value = '  +6.17%  '
$ws.Range('D14').Value = '2.589.75'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').Value = '''0.883'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').Value = '''14.39'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '42.923.34'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '0.0₃0997'
$ws.Range('E18').Value = '  +2.03%  '
$ws.Range('D19').Value = '''12.88'
$ws.Range('E19').Value = '  +2.93%  '
$ws.Range('D20').Value = '''6.63'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').Value = '''72.04'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').Value = '''254.65'
$ws.Range('E22').Value = '  -3.53%  '
$ws.Range('E23').Value = '  +1.20%  '
$ws.Range('D24').Value = '''2.12'
$ws.Range('E24').Value = '  -2.84%  '
$ws.Range('D25').Value = '''28.89'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '''10.24'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('D28').Value = '''37.61'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').Value = '''6.02'
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').Value = '''155.29'
$ws.Range('E31').Value = '  +2.56%  '
$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '''2.76'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '''2.16'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''3.39'
$ws.Range('E34').Value = '  -4.07%  '
$ws.Range('D35').Value = '''0.0805'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').Value = '''18.36'
$ws.Range('E36').Value = '  +11.64%  '
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = '''23.15'
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('D40').Value = '''3.42'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D44').Value = '2.075.95'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '''9.19'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('D47').Value = '''85.26'
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('D48').Value = '''76.63'
$ws.Range('E48').Value = '  +11.18%  '
$ws.Range('D49').Value = '''106.71'
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('D50').Value = '2.821.04'
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').Value = '''1.67'
